$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(1).ColumnWidth = 2.83203125
$ws.Columns.Item(1).BestFit = $true
